$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.500691
$ws.Range("H2").Value = 1.502073
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 56.98117766666667
$ws.Range("N2").Value = 170.943533
$ws.Range("O2").Value = 0.952030123851636
$ws.Range("P2").Value = 0.9520301238516359
$ws.Range("Q2").Value = 28.529962827101
$ws.Range("R2").Value = 256.769665443909
$ws.Range("S2").Value = 0.952030123851636
$ws.Range("T2").Value = 0.9520301238516359

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.500691
$ws.Range("H3").Value = 1.502073
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.516719
$ws.Range("N3").Value = 7.550157
$ws.Range("O3").Value = 0.04204883786863874
$ws.Range("P3").Value = 0.04204883786863874
$ws.Range("Q3").Value = 1.260098552829
$ws.Range("R3").Value = 11.340886975461
$ws.Range("S3").Value = 0.04204883786863874
$ws.Range("T3").Value = 0.04204883786863874

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.500691
$ws.Range("H4").Value = 1.502073
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.3543876666666666
$ws.Range("N4").Value = 1.063163
$ws.Range("O4").Value = 0.005921038279725251
$ws.Range("P4").Value = 0.005921038279725251
$ws.Range("Q4").Value = 0.177438715211
$ws.Range("R4").Value = 1.596948436899
$ws.Range("S4").Value = 0.005921038279725251
$ws.Range("T4").Value = 0.005921038279725251
